$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row()
if ($lastRow -lt 2) { $lastRow = 294 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value2()
    if ($v -eq 45203) {
        $cell.Value = 45204
    }
}
